$d = $word.ActiveDocument

# Remove the existing "_GoBack" bookmark (currently spans from the very
# start of the document to the end of the last paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Locate the end of " anuncios restantes" and append the new
# "(Acabar hoy 19/10/2017)" text as three separate runs, matching the
# run-splitting seen in the target document (the bookmark will sit
# between the date and the closing parenthesis).
$rng = $d.Content
$rng.Find.Execute("anuncios restantes", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertAfter("(Acabar hoy")
$rng.Collapse(0)
$rng.InsertAfter(" 19/10/2017")
$rng.Collapse(0)

# Remember this position (between the date and the closing parenthesis)
# for the "_GoBack" bookmark before inserting the final run, so the
# bookmark range does not end up sitting on a paragraph-end boundary.
$bmStart = $rng.Start
$rng.InsertAfter(")")

# Re-create the "_GoBack" bookmark as a zero-length bookmark right
# between the date run and the closing parenthesis run.
$bmRng = $d.Range($bmStart, $bmStart)
$d.Bookmarks.Add("_GoBack", $bmRng)
